$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column C ("max") - shifts D->C, E->D
$ws.Range("C1").EntireColumn.Delete()

# Delete rows 4 and 5 (even_MAG-GUT73847.fa, even_MAG-GUT74183.fa) - shifts rows 6-9 up to 4-7
$ws.Range("A4:A5").EntireRow.Delete()

# Now set updated values for remaining rows
$ws.Range("B2").Value = 2.800922724863074
$ws.Range("B3").Value = 1.570366708433458
$ws.Range("B4").Value = 1.575650301524327
$ws.Range("B5").Value = 3.353323903427377
$ws.Range("B6").Value = 1.631202863296863
$ws.Range("B7").Value = 2.83861884753728

$ws.Range("C2").Value = "s__CAG-988 sp003149915"
$ws.Range("C3").Value = "s__CAG-988 sp003149915"
$ws.Range("C4").Value = "s__CAG-988 sp003149915"
$ws.Range("C5").Value = "s__CAG-988 sp003149915"
$ws.Range("C6").Value = "s__CAG-988 sp003149915"
$ws.Range("C7").Value = "s__CAG-988 sp003149915"
